$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / non-numeric-looking values: direct assignment is safe ---
$ws.Range("D2").Value = "44.502.76"
$ws.Range("E2").Value = "  +3.91%  "
$ws.Range("D3").Value = "2.280.02"
$ws.Range("E3").Value = "  +2.46%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("E6").Value = "  +6.19%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("E10").Value = "  +3.63%  "
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("E14").Value = "  +2.58%  "
$ws.Range("D15").Value = "2.628.92"
$ws.Range("E15").Value = "  +2.56%  "
$ws.Range("E16").Value = "  +2.19%  "
$ws.Range("D17").Value = "2.279.97"
$ws.Range("E17").Value = "  +3.07%  "
$ws.Range("D18").Value = "44.356.36"
$ws.Range("E18").Value = "  +3.57%  "
$ws.Range("E19").Value = "  -5.85%  "
$ws.Range("E20").Value = "  +4.32%  "
$ws.Range("E21").Value = "  +1.46%  "
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("E24").Value = "  +1.13%  "
$ws.Range("E25").Value = "  +2.64%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  +1.29%  "
$ws.Range("E29").Value = "  +12.25%  "
$ws.Range("E30").Value = "  +2.21%  "
$ws.Range("E31").Value = "  +4.77%  "
$ws.Range("E32").Value = "  +0.70%  "
$ws.Range("E33").Value = "  -2.61%  "
$ws.Range("E34").Value = "  -0.71%  "
$ws.Range("E35").Value = "  +4.15%  "
$ws.Range("E36").Value = "  +10.28%  "
$ws.Range("E37").Value = "  -0.50%  "
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("E39").Value = "  +1.25%  "
$ws.Range("E40").Value = "  -0.23%  "
$ws.Range("E41").Value = "  +24.58%  "
$ws.Range("E42").Value = "  +0.80%  "
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").Value = "1.778.79"
$ws.Range("E44").Value = "  -8.59%  "
$ws.Range("E45").Value = "  -2.54%  "
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("E47").Value = "  +2.27%  "
$ws.Range("E48").Value = "  -0.39%  "
$ws.Range("E49").Value = "  -2.12%  "
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("E51").Value = "  +1.23%  "

# --- Numeric-looking strings: force text storage so Excel does not convert them to numbers ---
$numericTextCells = [ordered]@{
    "D5" = "320.61"
    "D6" = "106.16"
    "D9" = "0.572"
    "D10" = "38.78"
    "D11" = "0.0844"
    "D12" = "7.90"
    "D14" = "0.885"
    "D16" = "14.61"
    "D19" = "14.06"
    "D22" = "66.44"
    "D23" = "3.21"
    "D24" = "239.24"
    "D27" = "10.20"
    "D29" = "38.34"
    "D31" = "164.08"
    "D32" = "20.64"
    "D33" = "0.0886"
    "D35" = "2.04"
    "D37" = "3.19"
    "D40" = "4.46"
    "D41" = "15.61"
    "D45" = "86.89"
    "D48" = "60.36"
    "D49" = "74.84"
    "D50" = "1.70"
    "D51" = "8.69"
}
foreach ($cellRef in $numericTextCells.Keys) {
    $target = $ws.Range($cellRef)
    $target.NumberFormat = "@"
    $target.Value = $numericTextCells[$cellRef]
    $target.Style = "Normal"
}
